$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 81211
$ws.Range("C2").Value = 552.2850352784722
$ws.Range("D2").Value = 122.486354368838
$ws.Range("E2").Value = 407
$ws.Range("F2").Value = 469
$ws.Range("G2").Value = 514
$ws.Range("H2").Value = 590
$ws.Range("I2").Value = 1510

$ws.Range("B3").Value = 81211
$ws.Range("C3").Value = 45.75903485980964
$ws.Range("D3").Value = 4.404573265810024
$ws.Range("E3").Value = 32.43
$ws.Range("F3").Value = 42.41
$ws.Range("G3").Value = 45.66
$ws.Range("H3").Value = 48.89

$ws.Range("B4").Value = 81211
$ws.Range("C4").Value = 1.815599734026178
$ws.Range("D4").Value = 5.98815343961479
$ws.Range("F4").Value = 0.76
$ws.Range("G4").Value = 1.41
$ws.Range("H4").Value = 2.25
$ws.Range("I4").Value = 637.71

$ws.Range("B5").Value = 81211
$ws.Range("C5").Value = 321.6290285798722
$ws.Range("D5").Value = 3.794901112960051
$ws.Range("E5").Value = 304.31
$ws.Range("F5").Value = 319.24
$ws.Range("G5").Value = 321.56
$ws.Range("H5").Value = 324.74
$ws.Range("I5").Value = 330.38

$ws.Range("B6").Value = 81211
$ws.Range("C6").Value = 22.21511814901922
$ws.Range("D6").Value = 1.70474918211924
$ws.Range("E6").Value = 16.2
$ws.Range("F6").Value = 21.29
$ws.Range("G6").Value = 21.84
$ws.Range("H6").Value = 22.74
$ws.Range("I6").Value = 31.8

$ws.Range("B7").Value = 81211
$ws.Range("C7").Value = -76.28842151925232
$ws.Range("D7").Value = 24.07634987098373
$ws.Range("E7").Value = -122
$ws.Range("F7").Value = -96
$ws.Range("G7").Value = -72
$ws.Range("H7").Value = -56

$ws.Range("B8").Value = 81089
$ws.Range("C8").Value = 7.533825796347224
$ws.Range("D8").Value = 6.592628138518452
$ws.Range("E8").Value = -23.8
$ws.Range("F8").Value = 7.5
$ws.Range("I8").Value = 15

$ws.Range("B9").Value = 81211
$ws.Range("C9").Value = 9.327640344288335
$ws.Range("D9").Value = 1.685565108282402

$ws.Range("B10").Value = 81211
$ws.Range("C10").Value = 867.8316151752841
$ws.Range("D10").Value = 0.4604519342066086

$ws.Range("B11").Value = 81211
$ws.Range("C11").Value = 26773.71921291451
$ws.Range("D11").Value = 4334.889041327159
$ws.Range("E11").Value = 18017
$ws.Range("F11").Value = 23613.5
$ws.Range("G11").Value = 26795
$ws.Range("H11").Value = 30044.5
$ws.Range("I11").Value = 36719

$ws.Range("B12").Value = 81211
$ws.Range("C12").Value = 29901.93253377006
$ws.Range("D12").Value = 4766.061616018711
$ws.Range("E12").Value = 20280
$ws.Range("F12").Value = 26432
$ws.Range("G12").Value = 29924
$ws.Range("H12").Value = 33504
$ws.Range("I12").Value = 40030

$ws.Range("B13").Value = 81211
$ws.Range("C13").Value = 0.5573775832337985
$ws.Range("D13").Value = 0.5900888384686347

$ws.Range("B14").Value = 81211
$ws.Range("C14").Value = 23.91567644777185
$ws.Range("D14").Value = 13.42217716876898
$ws.Range("H14").Value = 39
$ws.Range("I14").Value = 43

$ws.Range("B15").Value = 81211
$ws.Range("C15").Value = 0.6728398862223098
$ws.Range("D15").Value = 0.7484391091290741

$ws.Range("B16").Value = 81211
$ws.Range("C16").Value = 1.835391757274261
$ws.Range("D16").Value = 1.673899526794093

$ws.Range("B17").Value = 81211
$ws.Range("C17").Value = 93.68842151925233
$ws.Range("D17").Value = 24.07634987097837
$ws.Range("F17").Value = 73.40000000000001
$ws.Range("G17").Value = 89.40000000000001
$ws.Range("H17").Value = 113.4
$ws.Range("I17").Value = 139.4

$ws.Range("B18").Value = 81089
$ws.Range("C18").Value = -85.32753334908405
$ws.Range("D18").Value = 21.726273662449
$ws.Range("E18").Value = -125.3773603942068
$ws.Range("F18").Value = -104.6286946522615
$ws.Range("G18").Value = -83.5149694202523
$ws.Range("H18").Value = -66.7376019773414

$ws.Range("B19").Value = 81089
$ws.Range("C19").Value = -77.79370755273685
$ws.Range("D19").Value = 26.14762289951132
$ws.Range("F19").Value = -96.41392685158225
$ws.Range("G19").Value = -72.41392685158225
$ws.Range("H19").Value = -56.2778545523916
$ws.Range("I19").Value = -33.49305820175223
